$wb = $excel.ActiveWorkbook

# --- Sheet 1: safety_orders ---
$ws1 = $wb.Worksheets.Item("safety_orders")

$ws1.Range("D2").Value = 1.778743315536
$ws1.Range("E2").Value = 1.818477856668
$ws1.Range("F2").Value = 1.83666263523468

$ws1.Range("D3").Value = 1.68486178223616
$ws1.Range("E3").Value = 1.75166981945208
$ws1.Range("F3").Value = 1.769186517646601
$ws1.Range("G3").Value = 4.766299910684991

$ws1.Range("D4").Value = 1.538406590250365
$ws1.Range("E4").Value = 1.645038204851223
$ws1.Range("F4").Value = 1.661488586899735

$ws1.Range("D5").Value = 1.309936490828614
$ws1.Range("E5").Value = 1.477487347839918
$ws1.Range("F5").Value = 1.492262221318318

$ws1.Range("D6").Value = 0.9535231356774201
$ws1.Range("E6").Value = 1.215505241758669
$ws1.Range("F6").Value = 1.227660294176256
$ws1.Range("G6").Value = 22.33005007975585

# --- Sheet 2: open_buy_orders ---
# Row 2 values change, row 3 (previously the second order) is removed entirely.
$ws2 = $wb.Worksheets.Item("open_buy_orders")

$ws2.Range("A2").Value = "O7QK74-SAHCK-DOHWBC"
$ws2.Range("B2").Value = 1.87679
$ws2.Rows.Item(3).Delete()

# --- Sheet 3: open_sell_orders ---
# Row 2 txid changes, and a new row 3 is appended.
$ws3 = $wb.Worksheets.Item("open_sell_orders")

$ws3.Range("A2").Value = "OC24JQ-P7CJG-5E45V2"
$ws3.Range("A3").Value = "OH5JAF-MVZKF-CEDQ2M"
